$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.844.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.365.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "658.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.62%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.423"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.10%  "

# Row 9
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("E10").Value = "  -6.27%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.363.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.209"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.522.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.08"
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000255"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.67%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.992.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.74%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.09%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.363.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.516"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -13.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "510.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.77%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.65%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.84%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.76%  "

# Row 30
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.141"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.185"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.86%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.558"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.73%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "524.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.97%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.151"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.16%  "

# Row 41
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.96%  "

# Row 43
$ws.Range("E43").Value = "  -1.30%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.853"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.75%  "

# Row 45
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0425"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.39%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.45%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.54%  "
